$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 81 with new values (previously "teste de jogo ps5" test entry,
# now replaced with the real "GTA Vice City" entry)
$ws.Range("A81").Value = "GTA Vice City"
$ws.Range("B81").Value = "em progresso"
$ws.Range("C81").Value = "PC"
$ws.Range("D81").Value = "Zerar"

# Remove row 82 entirely (the leftover "teste" entry), shifting dimension
# from A1:D82 down to A1:D81
$ws.Rows.Item(82).Delete()
